$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.291.18'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.932.17'
$ws.Range("E3").Value = '  +0.36%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = "'0.7502"
$ws.Range("E5").Value = '  +4.64%  '
$ws.Range("D6").Value = "'242.69"
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("D8").Value = "'27.83"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = "'0.3183"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = "'0.07160"
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("E11").Value = '  -1.03%  '
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("D13").Value = '1.919.55'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'5.396"
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = "'93.03"
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = "'14.55"
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").Value = '30.302.77'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = "'6.068"
$ws.Range("E18").Value = '  +5.32%  '
$ws.Range("D19").Value = "'251.53"
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").Value = "'0.000007963"
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").Value = '2.171.60'
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = "'6.676"
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("D25").Value = "'9.551"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").Value = "'164.64"
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").Value = "'19.10"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = "'0.1298"
$ws.Range("E28").Value = '  +2.51%  '
$ws.Range("D29").Value = "'2.193"
$ws.Range("E29").Value = '  -2.88%  '
$ws.Range("D30").Value = "'1.376"
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("D31").Value = "'1.545"
$ws.Range("D32").Value = "'4.417"
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = "'4.149"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "'0.05241"
$ws.Range("E34").Value = '  +2.16%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'1.324"
$ws.Range("E35").Value = '  +4.56%  '
$ws.Range("D36").Value = "'0.7565"
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("E37").Value = '  +0.89%  '
$ws.Range("D38").Value = "'0.01955"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = "'2.801"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = "'78.67"
$ws.Range("E40").Value = '  +1.30%  '
$ws.Range("D41").Value = "'6.511"
$ws.Range("E41").Value = '  +2.59%  '
$ws.Range("D42").Value = "'0.4527"
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").Value = "'1.979"
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = "'0.8415"
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = "'10.05"
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("D47").Value = "'7.694"
$ws.Range("E47").Value = '  +3.72%  '
$ws.Range("D48").Value = "'101.82"
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("D49").Value = "'37.77"
$ws.Range("E49").Value = '  +3.36%  '
$ws.Range("D50").Value = "'0.1216"
$ws.Range("E50").Value = '  +7.22%  '
$ws.Range("D51").Value = "'956.74"
$ws.Range("E51").Value = '  +3.49%  '
